$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header and data cells
$ws.Range("A1").Value = "Items"
$ws.Range("B1").Value = "Status"
$ws.Range("A2").Value = "testing,shopping,owrjdfnd,43545#@@,party!!"
$ws.Range("B2").Value = "complete,in_progress,complete,complete,in_progress"

# Update selection to match the target view
$ws.Range("B10").Select()
